$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# Replace every exact "+0" cell text with "0" (matches whole word so that
# "5+5" and similar strings are left untouched).
$find.Execute("+0", $true, $true, $false, $false, $false, $true, 1, $false, "0", 2)
